$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "BTC"
$ws.Range("C2").Value = "Bitcoin"
$ws.Range("D2").Value = 42549
$ws.Range("E2").Value = 833838111405
$ws.Range("F2").Value = 15784397386
$ws.Range("G2").Value = 1.15327

$ws.Range("B3").Value = "ETH"
$ws.Range("C3").Value = "Ethereum"
$ws.Range("D3").Value = 2312.96
$ws.Range("E3").Value = 278223852789
$ws.Range("F3").Value = 13770484760
$ws.Range("G3").Value = 0.66005

$ws.Range("B4").Value = "USDT"
$ws.Range("C4").Value = "Tether"
$ws.Range("D4").Value = 0.999816
$ws.Range("E4").Value = 91703323180
$ws.Range("F4").Value = 26972034979
$ws.Range("G4").Value = -0.06954

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "BNB"
$ws.Range("D5").Value = 317.71
$ws.Range("E5").Value = 48864314566
$ws.Range("F5").Value = 762319018
$ws.Range("G5").Value = -0.13952

$ws.Range("B6").Value = "SOL"
$ws.Range("C6").Value = "Solana"
$ws.Range("D6").Value = 104.85
$ws.Range("E6").Value = 45040650600
$ws.Range("F6").Value = 1801340041
$ws.Range("G6").Value = 1.67517

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "XRP"
$ws.Range("D7").Value = 0.623647
$ws.Range("E7").Value = 33763237785
$ws.Range("F7").Value = 587930520
$ws.Range("G7").Value = -0.03098

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "USDC"
$ws.Range("D8").Value = 0.9994729999999999
$ws.Range("E8").Value = 24624218935
$ws.Range("F8").Value = 7727557381
$ws.Range("G8").Value = -0.12624

$ws.Range("B9").Value = "STETH"
$ws.Range("C9").Value = "Lido Staked Ether"
$ws.Range("D9").Value = 2312.41
$ws.Range("E9").Value = 21303488758
$ws.Range("F9").Value = 14139218
$ws.Range("G9").Value = 0.75459

$ws.Range("B10").Value = "ADA"
$ws.Range("C10").Value = "Cardano"
$ws.Range("D10").Value = 0.606781
$ws.Range("E10").Value = 21251961193
$ws.Range("F10").Value = 366605693
$ws.Range("G10").Value = 0.67725

$ws.Range("B11").Value = "AVAX"
$ws.Range("C11").Value = "Avalanche"
$ws.Range("D11").Value = 40.15
$ws.Range("E11").Value = 14689094580
$ws.Range("F11").Value = 592484813
$ws.Range("G11").Value = 2.0011

$ws.Range("B12").Value = "DOGE"
$ws.Range("C12").Value = "Dogecoin"
$ws.Range("D12").Value = 0.09050800000000001
$ws.Range("E12").Value = 12902780919
$ws.Range("F12").Value = 280079440
$ws.Range("G12").Value = 0.25808

$ws.Range("B13").Value = "DOT"
$ws.Range("C13").Value = "Polkadot"
$ws.Range("D13").Value = 8.51
$ws.Range("E13").Value = 11191211124
$ws.Range("F13").Value = 281015881
$ws.Range("G13").Value = 3.16845

$ws.Range("B14").Value = "TRX"
$ws.Range("C14").Value = "TRON"
$ws.Range("D14").Value = 0.106548
$ws.Range("E14").Value = 9413708410
$ws.Range("F14").Value = 298570097
$ws.Range("G14").Value = 0.49362

$ws.Range("B15").Value = "MATIC"
$ws.Range("C15").Value = "Polygon"
$ws.Range("D15").Value = 0.983429
$ws.Range("E15").Value = 9110381862
$ws.Range("F15").Value = 359830884
$ws.Range("G15").Value = 2.17208

$ws.Range("B16").Value = "LINK"
$ws.Range("C16").Value = "Chainlink"
$ws.Range("D16").Value = 15.45
$ws.Range("E16").Value = 8592085839
$ws.Range("F16").Value = 336739296
$ws.Range("G16").Value = 1.47146

$ws.Range("B17").Value = "TON"
$ws.Range("C17").Value = "Toncoin"
$ws.Range("D17").Value = 2.28
$ws.Range("E17").Value = 7880132911
$ws.Range("F17").Value = 44025587
$ws.Range("G17").Value = -2.56502

$ws.Range("B18").Value = "WBTC"
$ws.Range("C18").Value = "Wrapped Bitcoin"
$ws.Range("D18").Value = 42530
$ws.Range("E18").Value = 6725943326
$ws.Range("F18").Value = 117952231
$ws.Range("G18").Value = 1.0643

$ws.Range("B19").Value = "ICP"
$ws.Range("C19").Value = "Internet Computer"
$ws.Range("D19").Value = 13.71
$ws.Range("E19").Value = 6264159141
$ws.Range("F19").Value = 805580925
$ws.Range("G19").Value = 36.85761

$ws.Range("B20").Value = "SHIB"
$ws.Range("C20").Value = "Shiba Inu"
$ws.Range("D20").Value = 0.00001056
$ws.Range("E20").Value = 6223907198
$ws.Range("F20").Value = 83318161
$ws.Range("G20").Value = 0.55955

$ws.Range("B21").Value = "UNI"
$ws.Range("C21").Value = "Uniswap"
$ws.Range("D21").Value = 7.66
$ws.Range("E21").Value = 5771585711
$ws.Range("F21").Value = 214019179
$ws.Range("G21").Value = 2.05048

$ws.Range("B22").Value = "LTC"
$ws.Range("C22").Value = "Litecoin"
$ws.Range("D22").Value = 73.88
$ws.Range("E22").Value = 5472217772
$ws.Range("F22").Value = 427636268
$ws.Range("G22").Value = 0.43561

$ws.Range("B23").Value = "DAI"
$ws.Range("C23").Value = "Dai"
$ws.Range("D23").Value = 0.9989130000000001
$ws.Range("E23").Value = 5236682697
$ws.Range("F23").Value = 211424426
$ws.Range("G23").Value = -0.03898

$ws.Range("B24").Value = "BCH"
$ws.Range("C24").Value = "Bitcoin Cash"
$ws.Range("D24").Value = 266.96
$ws.Range("E24").Value = 5231494920
$ws.Range("F24").Value = 300857479
$ws.Range("G24").Value = -3.40793

$ws.Range("B25").Value = "ATOM"
$ws.Range("C25").Value = "Cosmos Hub"
$ws.Range("D25").Value = 10.88
$ws.Range("E25").Value = 4146133046
$ws.Range("F25").Value = 178852410
$ws.Range("G25").Value = 0.84501

$ws.Range("B26").Value = "NEAR"
$ws.Range("C26").Value = "NEAR Protocol"
$ws.Range("D26").Value = 3.71
$ws.Range("E26").Value = 3768776847
$ws.Range("F26").Value = 212761672
$ws.Range("G26").Value = 0.04485

$ws.Range("B27").Value = "XLM"
$ws.Range("C27").Value = "Stellar"
$ws.Range("D27").Value = 0.131129
$ws.Range("E27").Value = 3704208892
$ws.Range("F27").Value = 95242832
$ws.Range("G27").Value = -3.51395

$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "LEO Token"
$ws.Range("D28").Value = 3.94
$ws.Range("E28").Value = 3650161960
$ws.Range("F28").Value = 1312282
$ws.Range("G28").Value = -0.71728

$ws.Range("B29").Value = "OP"
$ws.Range("C29").Value = "Optimism"
$ws.Range("D29").Value = 3.87
$ws.Range("E29").Value = 3523741775
$ws.Range("F29").Value = 316138700
$ws.Range("G29").Value = 5.91692

$ws.Range("B30").Value = "OKB"
$ws.Range("C30").Value = "OKB"
$ws.Range("D30").Value = 54.89
$ws.Range("E30").Value = 3289278678
$ws.Range("F30").Value = 4780484
$ws.Range("G30").Value = 0.53542

$ws.Range("B31").Value = "INJ"
$ws.Range("C31").Value = "Injective"
$ws.Range("D31").Value = 38.39
$ws.Range("E31").Value = 3234833395
$ws.Range("F31").Value = 233504804
$ws.Range("G31").Value = 4.77637

$ws.Range("B32").Value = "ETC"
$ws.Range("C32").Value = "Ethereum Classic"
$ws.Range("D32").Value = 22.39
$ws.Range("E32").Value = 3205672000
$ws.Range("F32").Value = 132182532
$ws.Range("G32").Value = -0.42793

$ws.Range("B33").Value = "FIL"
$ws.Range("C33").Value = "Filecoin"
$ws.Range("D33").Value = 6.38
$ws.Range("E33").Value = 3112565864
$ws.Range("F33").Value = 303018046
$ws.Range("G33").Value = 9.403180000000001

$ws.Range("B34").Value = "XMR"
$ws.Range("C34").Value = "Monero"
$ws.Range("D34").Value = 167.39
$ws.Range("E34").Value = 3037950711
$ws.Range("F34").Value = 63398344
$ws.Range("G34").Value = 2.27132

$ws.Range("B35").Value = "HBAR"
$ws.Range("C35").Value = "Hedera"
$ws.Range("D35").Value = 0.08870500000000001
$ws.Range("E35").Value = 2984347616
$ws.Range("F35").Value = 51347560
$ws.Range("G35").Value = 1.44135

$ws.Range("B36").Value = "APT"
$ws.Range("C36").Value = "Aptos"
$ws.Range("D36").Value = 9.59
$ws.Range("E36").Value = 2951402404
$ws.Range("F36").Value = 107270399
$ws.Range("G36").Value = -0.14932

$ws.Range("B37").Value = "IMX"
$ws.Range("C37").Value = "Immutable"
$ws.Range("D37").Value = 2.23
$ws.Range("E37").Value = 2945222558
$ws.Range("F37").Value = 65944518
$ws.Range("G37").Value = -1.04265

$ws.Range("B38").Value = "CRO"
$ws.Range("C38").Value = "Cronos"
$ws.Range("D38").Value = 0.099966
$ws.Range("E38").Value = 2645218835
$ws.Range("F38").Value = 7745997
$ws.Range("G38").Value = 0.91401

$ws.Range("B39").Value = "VET"
$ws.Range("C39").Value = "VeChain"
$ws.Range("D39").Value = 0.03503288
$ws.Range("E39").Value = 2545736406
$ws.Range("F39").Value = 79507761
$ws.Range("G39").Value = 0.71673

$ws.Range("B40").Value = "KAS"
$ws.Range("C40").Value = "Kaspa"
$ws.Range("D40").Value = 0.114377
$ws.Range("E40").Value = 2530611901
$ws.Range("F40").Value = 26881783
$ws.Range("G40").Value = 0.5282

$ws.Range("B41").Value = "LDO"
$ws.Range("C41").Value = "Lido DAO"
$ws.Range("D41").Value = 2.79
$ws.Range("E41").Value = 2469858099
$ws.Range("F41").Value = 64994509
$ws.Range("G41").Value = -5.4025

$ws.Range("B42").Value = "TUSD"
$ws.Range("C42").Value = "TrueUSD"
$ws.Range("D42").Value = 0.998
$ws.Range("E42").Value = 2305914612
$ws.Range("F42").Value = 176055088
$ws.Range("G42").Value = -0.29236

$ws.Range("B43").Value = "STX"
$ws.Range("C43").Value = "Stacks"
$ws.Range("D43").Value = 1.48
$ws.Range("E43").Value = 2112756990
$ws.Range("F43").Value = 72477977
$ws.Range("G43").Value = 6.7797

$ws.Range("B44").Value = "ARB"
$ws.Range("C44").Value = "Arbitrum"
$ws.Range("D44").Value = 1.64
$ws.Range("E44").Value = 2091112191
$ws.Range("F44").Value = 636633854
$ws.Range("G44").Value = 12.80495

$ws.Range("B45").Value = "QNT"
$ws.Range("C45").Value = "Quant"
$ws.Range("D45").Value = 142.99
$ws.Range("E45").Value = 2078948113
$ws.Range("F45").Value = 60813813
$ws.Range("G45").Value = 6.09684

$ws.Range("B46").Value = "MNT"
$ws.Range("C46").Value = "Mantle"
$ws.Range("D46").Value = 0.659605
$ws.Range("E46").Value = 2069351169
$ws.Range("F46").Value = 118036912
$ws.Range("G46").Value = 1.89642

$ws.Range("B47").Value = "BSV"
$ws.Range("C47").Value = "Bitcoin SV"
$ws.Range("D47").Value = 101.31
$ws.Range("E47").Value = 1955875065
$ws.Range("F47").Value = 319053720
$ws.Range("G47").Value = 4.97823

$ws.Range("B48").Value = "EGLD"
$ws.Range("C48").Value = "MultiversX"
$ws.Range("D48").Value = 69.78
$ws.Range("E48").Value = 1840986477
$ws.Range("F48").Value = 27600506
$ws.Range("G48").Value = 1.52936

$ws.Range("B49").Value = "TIA"
$ws.Range("C49").Value = "Celestia"
$ws.Range("D49").Value = 12.37
$ws.Range("E49").Value = 1807048999
$ws.Range("F49").Value = 79580098
$ws.Range("G49").Value = 3.1194

$ws.Range("B50").Value = "FDUSD"
$ws.Range("C50").Value = "First Digital USD"
$ws.Range("D50").Value = 1.001
$ws.Range("E50").Value = 1802969442
$ws.Range("F50").Value = 1835513314
$ws.Range("G50").Value = -0.21075

$ws.Range("B51").Value = "ALGO"
$ws.Range("C51").Value = "Algorand"
$ws.Range("D51").Value = 0.224498
$ws.Range("E51").Value = 1796537361
$ws.Range("F51").Value = 48908859
$ws.Range("G51").Value = 0.71004
